$wb = $excel.ActiveWorkbook
$wsUsers = $wb.Worksheets.Item("Users")
$wsEng = $wb.Worksheets.Item("Engagements")

# --- Engagements sheet: add new "EngagementLOB" column (header + existing rows) ---
$wsEng.Range("B1").Value = "EngagementLOB"
$wsEng.Range("B1").Font.Bold = $true
$wsEng.Range("B2").Value = "CF"

# --- Engagements sheet: new row for "Salem Harbor" ---
$wsEng.Range("A3").Value = "Salem Harbor"
$wsEng.Range("B3").Value = "FR"

# --- Users sheet: add two more team members ---
$wsUsers.Range("A3").Value = "Spencer Anderson"
$wsUsers.Range("A4").Value = "Karan Chopra"

# --- Engagements sheet: new row for "Arista Networks - Big Switch PPA" ---
$wsEng.Range("B4").Value = "FVA"
$wsEng.Range("A4").Value = "Arista Networks - Big Switch PPA"

$wsEng.Columns.Item(1).ColumnWidth = 27
$wsEng.Columns.Item(2).ColumnWidth = 14.166666666666666

$wsUsers.Range("A4").Select()
$wsEng.Range("A9").Select()
